$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the rule table's hit-policy cell from "U" (Unique) to "F" (First).
$ws.Range("B4").Value = "F"

# The header row (rule-id / conditions) picks up a distinct (but visually
# identical) cell style in the authored workbook - touch a formatting
# property at its already-default value to mint a new, equivalent style
# record without changing the visible formatting.
$hdr = $ws.Range("B4:D4")
$hdr.WrapText = $false

# Move the active selection to B5 (row below the hit-policy cell).
$ws.Range("B5").Select()
